# Auto-generated PowerShell COM-interop script
# Applies updated market-price / profit values to the Shinryu_Profits workbook
# as produced by the scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 633.3333
$ws.Range("I97").Value = 460
$ws.Range("J97").Value = 850
$ws.Range("K97").Value = 1380
$ws.Range("L97").Value = 2550
$ws.Range("M97").Value = -884
$ws.Range("N97").Value = -3542
$ws.Range("H137").Value = 44924.824
$ws.Range("I137").Value = 870.7143
$ws.Range("J137").Value = 113453.445
$ws.Range("K137").Value = 2612.1429
$ws.Range("L137").Value = 340360.335
$ws.Range("M137").Value = -62.14289999999983
$ws.Range("N137").Value = -345460.335
$ws.Range("H138").Value = 4142.609
$ws.Range("I138").Value = 1512.25
$ws.Range("J138").Value = 4696.3687
$ws.Range("K138").Value = 4536.75
$ws.Range("L138").Value = 14089.1061
$ws.Range("M138").Value = 603.25
$ws.Range("N138").Value = -24369.1061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2356.5
$ws.Range("I61").Value = 1169.0605
$ws.Range("J61").Value = 6710.4443
$ws.Range("K61").Value = 1169.0605
$ws.Range("L61").Value = 6710.4443
$ws.Range("M61").Value = -957.0605
$ws.Range("N61").Value = -7134.4443
$ws.Range("H74").Value = 4900.2964
$ws.Range("I74").Value = 5336.091
$ws.Range("J74").Value = 2982.8
$ws.Range("K74").Value = 5336.091
$ws.Range("L74").Value = 2982.8
$ws.Range("M74").Value = -4462.091
$ws.Range("N74").Value = -4730.8
$ws.Range("H77").Value = 4900.2964
$ws.Range("I77").Value = 5336.091
$ws.Range("J77").Value = 2982.8
$ws.Range("K77").Value = 26680.455
$ws.Range("L77").Value = 14914
$ws.Range("M77").Value = -22312.455
$ws.Range("N77").Value = -23650
$ws.Range("H88").Value = 1863.25
$ws.Range("J88").Value = 1733.3334
$ws.Range("L88").Value = 1733.3334
$ws.Range("N88").Value = -2545.3334
$ws.Range("H91").Value = 1863.25
$ws.Range("J91").Value = 1733.3334
$ws.Range("L91").Value = 1733.3334
$ws.Range("N91").Value = -4541.3334
$ws.Range("H122").Value = 1790.8182
$ws.Range("I122").Value = 1790.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5372.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2922.4546
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2527.7693
$ws.Range("I132").Value = 1781
$ws.Range("J132").Value = 3722.6
$ws.Range("K132").Value = 5343
$ws.Range("L132").Value = 11167.8
$ws.Range("M132").Value = -2813
$ws.Range("N132").Value = -16227.8
$ws.Range("H136").Value = 2356.5
$ws.Range("I136").Value = 1169.0605
$ws.Range("J136").Value = 6710.4443
$ws.Range("K136").Value = 3507.1815
$ws.Range("L136").Value = 20131.3329
$ws.Range("M136").Value = -957.1815000000001
$ws.Range("N136").Value = -25231.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 47351.777
$ws.Range("J92").Value = 47351.777
$ws.Range("L92").Value = 47351.777
$ws.Range("N92").Value = -52343.777
$ws.Range("H105").Value = 2776.72
$ws.Range("I105").Value = 1641.25
$ws.Range("J105").Value = 2875.4565
$ws.Range("K105").Value = 1641.25
$ws.Range("L105").Value = 2875.4565
$ws.Range("M105").Value = 105.75
$ws.Range("N105").Value = -6369.4565
$ws.Range("H134").Value = 1734.9131
$ws.Range("I134").Value = 1540.409
$ws.Range("J134").Value = 6014
$ws.Range("K134").Value = 4621.227000000001
$ws.Range("L134").Value = 18042
$ws.Range("M134").Value = -2086.227000000001
$ws.Range("N134").Value = -23112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2395.027
$ws.Range("I31").Value = 1552.4667
$ws.Range("K31").Value = 1552.4667
$ws.Range("M31").Value = -1257.4667
$ws.Range("H34").Value = 2395.027
$ws.Range("I34").Value = 1552.4667
$ws.Range("K34").Value = 1552.4667
$ws.Range("M34").Value = -1350.4667
$ws.Range("H58").Value = 1620.8846
$ws.Range("I58").Value = 1310.7
$ws.Range("J58").Value = 2654.8333
$ws.Range("K58").Value = 1310.7
$ws.Range("L58").Value = 2654.8333
$ws.Range("M58").Value = -1107.7
$ws.Range("N58").Value = -3060.8333
$ws.Range("H132").Value = 4061.1
$ws.Range("I132").Value = 2366.3333
$ws.Range("K132").Value = 7098.999899999999
$ws.Range("M132").Value = -4568.999899999999
$ws.Range("H134").Value = 2475.1304
$ws.Range("I134").Value = 1434
$ws.Range("J134").Value = 5425
$ws.Range("K134").Value = 4302
$ws.Range("L134").Value = 16275
$ws.Range("M134").Value = -1767
$ws.Range("N134").Value = -21345
$ws.Range("H136").Value = 1620.8846
$ws.Range("I136").Value = 1310.7
$ws.Range("J136").Value = 2654.8333
$ws.Range("K136").Value = 3932.1
$ws.Range("L136").Value = 7964.499899999999
$ws.Range("M136").Value = -1382.1
$ws.Range("N136").Value = -13064.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 784.25
$ws.Range("I5").Value = 753.4286
$ws.Range("K5").Value = 2260.2858
$ws.Range("M5").Value = -2148.2858
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("K56").Value = 15000
$ws.Range("M56").Value = -14470
$ws.Range("H122").Value = 6612.75
$ws.Range("J122").Value = 1378.2222
$ws.Range("L122").Value = 12403.9998
$ws.Range("N122").Value = -17303.9998
$ws.Range("H135").Value = 784.25
$ws.Range("I135").Value = 753.4286
$ws.Range("K135").Value = 6780.8574
$ws.Range("M135").Value = -4245.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4493.3335
$ws.Range("I132").Value = 4638.1177
$ws.Range("K132").Value = 13914.3531
$ws.Range("M132").Value = -11384.3531
$ws.Range("H136").Value = 1269377.9
$ws.Range("J136").Value = 1269377.9
$ws.Range("L136").Value = 3808133.7
$ws.Range("N136").Value = -3813233.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 36142.855
$ws.Range("J87").Value = 36142.855
$ws.Range("L87").Value = 36142.855
$ws.Range("N87").Value = -38388.855
$ws.Range("H88").Value = 27490
$ws.Range("J88").Value = 27490
$ws.Range("L88").Value = 27490
$ws.Range("N88").Value = -28346
$ws.Range("H90").Value = 36142.855
$ws.Range("J90").Value = 36142.855
$ws.Range("L90").Value = 108428.565
$ws.Range("N90").Value = -119660.565
$ws.Range("H91").Value = 27490
$ws.Range("J91").Value = 27490
$ws.Range("L91").Value = 27490
$ws.Range("N91").Value = -30454
$ws.Range("H132").Value = 13775.692
$ws.Range("I132").Value = 19399.334
$ws.Range("J132").Value = 6107.091
$ws.Range("K132").Value = 58198.00199999999
$ws.Range("L132").Value = 18321.273
$ws.Range("M132").Value = -55668.00199999999
$ws.Range("N132").Value = -23381.273
$ws.Range("H136").Value = 1737.721
$ws.Range("I136").Value = 1332.7567
$ws.Range("K136").Value = 3998.2701
$ws.Range("M136").Value = -1448.2701

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 35326.668
$ws.Range("J114").Value = 35326.668
$ws.Range("L114").Value = 35326.668
$ws.Range("N114").Value = -44004.668
$ws.Range("H132").Value = 3131.55
$ws.Range("I132").Value = 2757.5557
$ws.Range("J132").Value = 3437.5454
$ws.Range("K132").Value = 8272.667099999999
$ws.Range("L132").Value = 10312.6362
$ws.Range("M132").Value = -5742.667099999999
$ws.Range("N132").Value = -15372.6362
$ws.Range("H136").Value = 4149.5405
$ws.Range("I136").Value = 4129.9355
$ws.Range("K136").Value = 12389.8065
$ws.Range("M136").Value = -9839.806499999999
